$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 83; this shifts the existing rows 83-153
# down to 84-154 (and the sheet dimension grows from A1:T153 to A1:T154).
$ws.Rows(83).Insert()

# Populate the newly inserted row 83 with its data.
$ws.Cells.Item(83, 1).Value = 8
$ws.Cells.Item(83, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(83, 3).Value = "Coquimbo"
$ws.Cells.Item(83, 4).Value = 45062
$ws.Cells.Item(83, 5).Value = 4
$ws.Cells.Item(83, 6).Value = "Fruta"
$ws.Cells.Item(83, 7).Value = 100109
$ws.Cells.Item(83, 8).Value = "Uva"
$ws.Cells.Item(83, 9).Value = 100109001
$ws.Cells.Item(83, 10).Value = "Uva"
$ws.Cells.Item(83, 11).Value = "Red Globe"
$ws.Cells.Item(83, 12).Value = "Primera"
$ws.Cells.Item(83, 13).Value = 400
$ws.Cells.Item(83, 14).Value = 9000
$ws.Cells.Item(83, 15).Value = 10000
$ws.Cells.Item(83, 16).Value = 9500
$ws.Cells.Item(83, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(83, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(83, 19).Value = 528
$ws.Cells.Item(83, 20).Value = 18
